# Updates cached market-board calc columns (H:N) across several Leve
# sheets, per the scheduled-runner refresh. A couple of rows also lose
# their HQ/NQ profit cell entirely when the corresponding price is 0
# (matches the existing convention used throughout these sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 94
$ws.Range("H94").Value = 13624.167
$ws.Range("I94").Value = 2898
$ws.Range("J94").Value = 21285.715
$ws.Range("K94").Value = 2898
$ws.Range("L94").Value = 21285.715
$ws.Range("M94").Value = -2447
$ws.Range("N94").Value = -22187.715

# Row 99
$ws.Range("H99").Value = 217
$ws.Range("I99").Value = 217
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 651
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 847
$ws.Range("N99").ClearContents()

# Row 100
$ws.Range("H100").Value = 62502164
$ws.Range("I100").Value = 2500.8333
$ws.Range("J100").Value = 250001150
$ws.Range("K100").Value = 2500.8333
$ws.Range("L100").Value = 250001150
$ws.Range("M100").Value = -1959.8333
$ws.Range("N100").Value = -250002232

# Row 125
$ws.Range("H125").Value = 1252.3572
$ws.Range("I125").Value = 698
$ws.Range("J125").Value = 1474.1
$ws.Range("K125").Value = 6282
$ws.Range("L125").Value = 13266.9
$ws.Range("M125").Value = -3822
$ws.Range("N125").Value = -18186.9

# Row 129
$ws.Range("H129").Value = 1219.5156
$ws.Range("I129").Value = 567.2
$ws.Range("J129").Value = 1340.3148
$ws.Range("K129").Value = 1701.6
$ws.Range("L129").Value = 4020.9444
$ws.Range("M129").Value = 3298.4
$ws.Range("N129").Value = -14020.9444

# Row 141
$ws.Range("H141").Value = 2007.1875
$ws.Range("I141").Value = 1447.3846
$ws.Range("J141").Value = 4433
$ws.Range("K141").Value = 4342.1538
$ws.Range("L141").Value = 13299
$ws.Range("M141").Value = 837.8462
$ws.Range("N141").Value = -23659

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1447.4193
$ws.Range("I61").Value = 1644.9565
$ws.Range("J61").Value = 879.5
$ws.Range("K61").Value = 1644.9565
$ws.Range("L61").Value = 879.5
$ws.Range("M61").Value = -1432.9565
$ws.Range("N61").Value = -1303.5

# Row 97
$ws.Range("H97").Value = 2482.88
$ws.Range("I97").Value = 1437
$ws.Range("J97").Value = 4051.7
$ws.Range("K97").Value = 1437
$ws.Range("L97").Value = 4051.7
$ws.Range("M97").Value = -941
$ws.Range("N97").Value = -5043.7

# Row 136
$ws.Range("H136").Value = 1447.4193
$ws.Range("I136").Value = 1644.9565
$ws.Range("J136").Value = 879.5
$ws.Range("K136").Value = 4934.8695
$ws.Range("L136").Value = 2638.5
$ws.Range("M136").Value = -2384.8695
$ws.Range("N136").Value = -7738.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 688.3946999999999
$ws.Range("I94").Value = 602.1
$ws.Range("J94").Value = 1012
$ws.Range("K94").Value = 602.1
$ws.Range("L94").Value = 1012
$ws.Range("M94").Value = -151.1
$ws.Range("N94").Value = -1914

# Row 134
$ws.Range("H134").Value = 132323.56
$ws.Range("I134").Value = 151850.45
$ws.Range("J134").Value = 2144.3333
$ws.Range("K134").Value = 455551.35
$ws.Range("L134").Value = 6432.999899999999
$ws.Range("M134").Value = -453016.35

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1049.6
$ws.Range("I58").Value = 1008.05
$ws.Range("J58").Value = 1132.7
$ws.Range("K58").Value = 1008.05
$ws.Range("L58").Value = 1132.7
$ws.Range("M58").Value = -805.05
$ws.Range("N58").Value = -1538.7

# Row 99
$ws.Range("H99").Value = 3120.4546
$ws.Range("I99").Value = 4662.4
$ws.Range("J99").Value = 1835.5
$ws.Range("K99").Value = 4662.4
$ws.Range("L99").Value = 1835.5
$ws.Range("M99").Value = -3164.4
$ws.Range("N99").Value = -4831.5

# Row 105
$ws.Range("H105").Value = 1875.55
$ws.Range("I105").Value = 738.1818
$ws.Range("J105").Value = 3265.6667
$ws.Range("K105").Value = 738.1818
$ws.Range("L105").Value = 3265.6667
$ws.Range("M105").Value = 1008.8182
$ws.Range("N105").Value = -6759.6667

# Row 126
$ws.Range("H126").Value = 3120.4546
$ws.Range("I126").Value = 4662.4
$ws.Range("J126").Value = 1835.5
$ws.Range("K126").Value = 13987.2
$ws.Range("L126").Value = 5506.5
$ws.Range("M126").Value = -11517.2
$ws.Range("N126").Value = -10446.5

# Row 132
$ws.Range("H132").Value = 2636.9375
$ws.Range("I132").Value = 2412.8
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 7238.400000000001
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -4708.400000000001
$ws.Range("N132").Value = -23057

# Row 136
$ws.Range("H136").Value = 1049.6
$ws.Range("I136").Value = 1008.05
$ws.Range("J136").Value = 1132.7
$ws.Range("K136").Value = 3024.15
$ws.Range("L136").Value = 3398.1
$ws.Range("M136").Value = -474.1499999999996
$ws.Range("N136").Value = -8498.1

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 52742.543
$ws.Range("I122").Value = 382.77777
$ws.Range("J122").Value = 59287.516
$ws.Range("K122").Value = 3444.99993
$ws.Range("L122").Value = 533587.6440000001
$ws.Range("M122").Value = -994.9999299999999
$ws.Range("N122").Value = -538487.6440000001

# Row 125
$ws.Range("H125").Value = 2416.6667
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 2545.4546
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 7636.3638
$ws.Range("M125").Value = 1920
$ws.Range("N125").Value = -17476.3638

# Row 126
$ws.Range("H126").Value = 2095.238
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2095.238
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 6285.714
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -16165.714

# Row 132
$ws.Range("H132").Value = 1013538
$ws.Range("I132").Value = 1645624.2
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 14810617.8
$ws.Range("L132").Value = 19800
$ws.Range("M132").Value = -14808087.8
$ws.Range("N132").Value = -24860

# Row 133
$ws.Range("H133").Value = 7671.4287
$ws.Range("I133").Value = 5750
$ws.Range("J133").Value = 8440
$ws.Range("K133").Value = 17250
$ws.Range("L133").Value = 25320
$ws.Range("M133").Value = -12190
$ws.Range("N133").Value = -35440

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2704.3845
$ws.Range("I122").Value = 2724.9
$ws.Range("J122").Value = 2636
$ws.Range("K122").Value = 8174.700000000001
$ws.Range("L122").Value = 7908
$ws.Range("M122").Value = -5724.700000000001
$ws.Range("N122").Value = -12808

# Row 126
$ws.Range("H126").Value = 4853.3335
$ws.Range("I126").Value = 2272.7273
$ws.Range("J126").Value = 7036.923
$ws.Range("K126").Value = 6818.1819
$ws.Range("L126").Value = 21110.769
$ws.Range("M126").Value = -4348.1819
$ws.Range("N126").Value = -26050.769

# Row 132
$ws.Range("H132").Value = 4497.5415
$ws.Range("I132").Value = 4394.0586
$ws.Range("J132").Value = 4748.857
$ws.Range("K132").Value = 13182.1758
$ws.Range("L132").Value = 14246.571
$ws.Range("M132").Value = -10652.1758
$ws.Range("N132").Value = -19306.571

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2797.8
$ws.Range("I40").Value = 2333
$ws.Range("J40").Value = 3495
$ws.Range("K40").Value = 2333
$ws.Range("L40").Value = 3495
$ws.Range("M40").Value = -2197
$ws.Range("N40").Value = -3767

# Row 92
$ws.Range("H92").Value = 31555.4
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 31555.4
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 31555.4
$ws.Range("N92").Value = -36547.4

# Row 93
$ws.Range("H93").Value = 1907.0834
$ws.Range("I93").Value = 1653.4445
$ws.Range("J93").Value = 2668
$ws.Range("K93").Value = 1653.4445
$ws.Range("L93").Value = 2668
$ws.Range("M93").Value = -405.4445000000001
$ws.Range("N93").Value = -5164

# Row 122
$ws.Range("H122").Value = 2967.7827
$ws.Range("I122").Value = 2208.25
$ws.Range("J122").Value = 3796.3635
$ws.Range("K122").Value = 6624.75
$ws.Range("L122").Value = 11389.0905
$ws.Range("M122").Value = -4174.75
$ws.Range("N122").Value = -16289.0905

# Row 136
$ws.Range("H136").Value = 1983.8572
$ws.Range("I136").Value = 1847.8636
$ws.Range("J136").Value = 2482.5
$ws.Range("K136").Value = 5543.5908
$ws.Range("L136").Value = 7447.5
$ws.Range("M136").Value = -2993.5908

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1650
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 1300
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 1300
$ws.Range("M96").Value = -627
$ws.Range("N96").Value = -4046

# Row 107
$ws.Range("H107").Value = 1568
$ws.Range("I107").Value = 1418.375
$ws.Range("J107").Value = 1739
$ws.Range("K107").Value = 4255.125
$ws.Range("L107").Value = 5217
$ws.Range("M107").Value = -2335.125
$ws.Range("N107").Value = -9057

# Row 122
$ws.Range("H122").Value = 4497.3125
$ws.Range("I122").Value = 6300.7144
$ws.Range("J122").Value = 3094.6667
$ws.Range("K122").Value = 18902.1432
$ws.Range("L122").Value = 9284.000100000001
$ws.Range("M122").Value = -16452.1432
$ws.Range("N122").Value = -14184.0001

# Row 136
$ws.Range("H136").Value = 2102.7827
$ws.Range("I136").Value = 1829.2
$ws.Range("J136").Value = 3926.6667
$ws.Range("K136").Value = 5487.6
$ws.Range("L136").Value = 11780.0001
$ws.Range("M136").Value = -2937.6
$ws.Range("N136").Value = -16880.0001
